$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.2212765957446808
$ws.Range("C2").Value = 0.4851063829787234
$ws.Range("J2").Value = 0.01702127659574468
$ws.Range("P2").Value = 0.1531914893617021
$ws.Range("S2").Value = 0.1234042553191489

# Row 3
$ws.Range("B3").Value = 0.008695652173913044
$ws.Range("J3").Value = 0.05217391304347826
$ws.Range("P3").Value = 0.6608695652173913
$ws.Range("S3").Value = 0.2782608695652174

# Row 4
$ws.Range("J4").Value = 0.02083333333333333
$ws.Range("P4").Value = 0.6875
$ws.Range("S4").Value = 0.2916666666666667

# Row 6
$ws.Range("B6").Value = 0.04310344827586207
$ws.Range("D6").Value = 0.01724137931034483
$ws.Range("F6").Value = 0.03017241379310345
$ws.Range("J6").Value = 0.3663793103448276
$ws.Range("O6").Value = 0.008620689655172414
$ws.Range("Q6").Value = 0.1810344827586207
$ws.Range("R6").Value = 0.05603448275862069
$ws.Range("S6").Value = 0.2974137931034483

# Row 7
$ws.Range("B7").Value = 0.124223602484472
$ws.Range("D7").Value = 0.0124223602484472
$ws.Range("F7").Value = 0.06211180124223602
$ws.Range("J7").Value = 0.1055900621118012
$ws.Range("O7").Value = 0.006211180124223602
$ws.Range("Q7").Value = 0.1925465838509317
$ws.Range("R7").Value = 0.09937888198757763
$ws.Range("S7").Value = 0.3975155279503105

# Row 8
$ws.Range("B8").Value = 0.05856832971800434
$ws.Range("D8").Value = 0.02386117136659436
$ws.Range("F8").Value = 0.06941431670281996
$ws.Range("J8").Value = 0.1214750542299349
$ws.Range("O8").Value = 0.01518438177874186
$ws.Range("Q8").Value = 0.1670281995661605
$ws.Range("R8").Value = 0.1540130151843818
$ws.Range("S8").Value = 0.3904555314533623

# Row 9
$ws.Range("B9").Value = 0.05696202531645569
$ws.Range("D9").Value = 0.0189873417721519
$ws.Range("F9").Value = 0.06962025316455696
$ws.Range("J9").Value = 0.08860759493670886
$ws.Range("O9").Value = 0.02531645569620253
$ws.Range("Q9").Value = 0.2151898734177215
$ws.Range("R9").Value = 0.1392405063291139
$ws.Range("S9").Value = 0.3860759493670886

# Row 10
$ws.Range("B10").Value = 0.09024979854955681
$ws.Range("D10").Value = 0.0225624496373892
$ws.Range("E10").Value = 0.0008058017727639
$ws.Range("F10").Value = 0.064464141821112
$ws.Range("J10").Value = 0.1095890410958904
$ws.Range("O10").Value = 0.016116035455278
$ws.Range("Q10").Value = 0.1909750201450443
$ws.Range("R10").Value = 0.1232876712328767
$ws.Range("S10").Value = 0.3819500402900886

# Row 11
$ws.Range("G11").Value = 0.1388888888888889
$ws.Range("J11").Value = 0.1111111111111111
$ws.Range("K11").Value = 0.1904761904761905
$ws.Range("L11").Value = 0.5555555555555556
$ws.Range("S11").Value = 0.003968253968253968

# Row 12
$ws.Range("G12").Value = 0.7375886524822695
$ws.Range("J12").Value = 0.1843971631205674
$ws.Range("K12").Value = 0.007092198581560284
$ws.Range("L12").Value = 0.0425531914893617
$ws.Range("S12").Value = 0.02836879432624113

# Row 13
$ws.Range("G13").Value = 0.7027027027027027
$ws.Range("J13").Value = 0.2972972972972973

# Row 15
$ws.Range("F15").Value = 0.02727272727272727
$ws.Range("H15").Value = 0.2
$ws.Range("I15").Value = 0.03636363636363636
$ws.Range("J15").Value = 0.4
$ws.Range("K15").Value = 0.05454545454545454
$ws.Range("M15").Value = 0.03181818181818181
$ws.Range("O15").Value = 0.03181818181818181
$ws.Range("S15").Value = 0.2181818181818182

# Row 16
$ws.Range("F16").Value = 0.03546099290780142
$ws.Range("H16").Value = 0.2340425531914894
$ws.Range("I16").Value = 0.0851063829787234
$ws.Range("J16").Value = 0.3333333333333333
$ws.Range("K16").Value = 0.1276595744680851
$ws.Range("M16").Value = 0.02127659574468085
$ws.Range("O16").Value = 0.03546099290780142
$ws.Range("S16").Value = 0.1276595744680851

# Row 17
$ws.Range("F17").Value = 0.02857142857142857
$ws.Range("H17").Value = 0.1904761904761905
$ws.Range("I17").Value = 0.05476190476190476
$ws.Range("J17").Value = 0.4595238095238095
$ws.Range("K17").Value = 0.1119047619047619
$ws.Range("M17").Value = 0.004761904761904762
$ws.Range("O17").Value = 0.07380952380952381
$ws.Range("S17").Value = 0.0761904761904762

# Row 18
$ws.Range("F18").Value = 0.05454545454545454
$ws.Range("H18").Value = 0.1527272727272727
$ws.Range("I18").Value = 0.05818181818181818
$ws.Range("J18").Value = 0.4872727272727272
$ws.Range("K18").Value = 0.08727272727272728
$ws.Range("M18").Value = 0.01090909090909091
$ws.Range("O18").Value = 0.06909090909090909
$ws.Range("S18").Value = 0.08

# Row 19
$ws.Range("F19").Value = 0.03025064822817632
$ws.Range("H19").Value = 0.22990492653414
$ws.Range("I19").Value = 0.08556611927398444
$ws.Range("J19").Value = 0.3621434745030251
$ws.Range("K19").Value = 0.08297320656871218
$ws.Range("M19").Value = 0.01987899740708729
$ws.Range("N19").Value = 0.001728608470181504
$ws.Range("O19").Value = 0.07692307692307693
$ws.Range("S19").Value = 0.1106309420916163
